# Updates the cryptos price list (Coin/Link/Price/Volume(1h)) on the active
# worksheet to match the latest scrape, per the Sep 30 2024 GitHub Actions run.
#
# Column D ('Price') holds figures stored as *text* in the workbook (not
# numbers), so values like '65.368.54' or '0.0₆0241' round-trip exactly.
# Assigning a numeric-looking string via COM .Value auto-coerces it to a real
# number (normal Excel behaviour), so those assignments are prefixed with a
# leading apostrophe to force a text entry, matching the source cell's text type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.637.14"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").Value = "2.613.97"
$ws.Range("E3").Value = "  -2.37%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'578.66"
$ws.Range("E5").Value = "  -3.82%  "
$ws.Range("D6").Value = "'157.25"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("D7").Value = "'0.650"
$ws.Range("E7").Value = "  +5.98%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.122"
$ws.Range("E9").Value = "  -5.64%  "
$ws.Range("D10").Value = "'5.79"
$ws.Range("E10").Value = "  -0.95%  "
$ws.Range("D11").Value = "'0.389"
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'28.38"
$ws.Range("E13").Value = "  -2.87%  "
$ws.Range("D14").Value = "'0.0000185"
$ws.Range("E14").Value = "  -7.61%  "
$ws.Range("D15").Value = "3.089.34"
$ws.Range("E15").Value = "  -2.12%  "
$ws.Range("D16").Value = "64.690.80"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "2.641.32"
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "'12.24"
$ws.Range("E18").Value = "  -4.30%  "
$ws.Range("D19").Value = "'4.67"
$ws.Range("E19").Value = "  -2.86%  "
$ws.Range("D20").Value = "'347.41"
$ws.Range("E20").Value = "  -1.57%  "
$ws.Range("D21").Value = "'7.32"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'68.33"
$ws.Range("E23").Value = "  -1.89%  "
$ws.Range("D24").Value = "'1.71"
$ws.Range("E24").Value = "  +2.20%  "
$ws.Range("D25").Value = "'0.0000109"
$ws.Range("E25").Value = "  -7.40%  "
$ws.Range("D26").Value = "'9.35"
$ws.Range("E26").Value = "  -2.96%  "
$ws.Range("D27").Value = "'1.59"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "'0.162"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'540.81"
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D31").Value = "'7.95"
$ws.Range("E31").Value = "  -1.10%  "
$ws.Range("D32").Value = "'2.09"
$ws.Range("E32").Value = "  -1.79%  "
$ws.Range("D33").Value = "'1.74"
$ws.Range("E33").Value = "  -1.43%  "
$ws.Range("D34").Value = "'6.40"
$ws.Range("E34").Value = "  -2.11%  "
$ws.Range("D35").Value = "'5.37"
$ws.Range("E35").Value = "  -2.11%  "
$ws.Range("D36").Value = "'0.412"
$ws.Range("E36").Value = "  -2.60%  "
$ws.Range("D37").Value = "'20.01"
$ws.Range("E37").Value = "  -3.13%  "
$ws.Range("D38").Value = "'1.00"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("D39").Value = "'1.93"
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "'150.05"
$ws.Range("E40").Value = "  -4.99%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'158.82"
$ws.Range("E42").Value = "  -2.64%  "
$ws.Range("B43").Value = "dogwifhat"
$ws.Range("C43").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D43").Value = "'2.41"
$ws.Range("E43").Value = "  +1.55%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'3.99"
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").Value = "'0.0606"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'22.67"
$ws.Range("E46").Value = "  -0.73%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.634"
$ws.Range("E47").Value = "  -1.44%  "
$ws.Range("E48").Value = "  +2.72%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0250"
$ws.Range("E49").Value = "  -3.09%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'19.16"
$ws.Range("E50").Value = "  -5.03%  "
$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "0.0₆0239"
$ws.Range("E51").Value = "  -6.49%  "
